# Updated CVDs for the month
# Applies the cell-value changes described by the OOXML diff to the six
# affected worksheets: Bangkrang Nonthaburi, Yueyang China, Ciserano Italy,
# Changzhou Epc China, Jiaxing China and Suzhou China.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a set of Cell -> Value pairs onto a worksheet.
# ---------------------------------------------------------------------------
function Set-Cells {
    param(
        $ws,
        $pairs
    )
    foreach ($ref in $pairs.Keys) {
        $ws.Range($ref).Value = $pairs[$ref]
    }
}

# ---------------------------------------------------------------------------
# 1) Bangkrang Nonthaburi
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Bangkrang Nonthaburi")

Set-Cells $ws1 @{
    "E2" = 0.0754; "E3" = 0.0754; "E4" = 0.0754
    "O4" = 0.0361; "P4" = 0; "Q4" = 0; "R4" = 0; "S4" = 0; "T4" = 0; "U4" = 0; "V4" = 0; "W4" = 0

    "E5" = 0.333333333333333; "E6" = 0.333333333333333; "E7" = 0.333333333333333
    "O7" = 0
    "P7" = 0.333333333333333; "Q7" = 0.333333333333333; "R7" = 0.333333333333333
    "S7" = 0.333333333333333; "T7" = 0.333333333333333; "U7" = 0.333333333333333
    "V7" = 0.333333333333333; "W7" = 0.333333333333333

    "E8" = 0.0809; "E9" = 0.0809; "E10" = 0.0809
    "O10" = 0.0113
    "P10" = 0.0115583333333333; "Q10" = 0.0115583333333333; "R10" = 0.034675
    "S10" = 0.0115583333333333; "T10" = 0.0115583333333333; "U10" = 0.0115583333333333
    "V10" = 0.034675; "W10" = 0.1387
}

# ---------------------------------------------------------------------------
# 2) Yueyang China
# ---------------------------------------------------------------------------
$ws11 = $wb.Worksheets.Item("Yueyang China")

Set-Cells $ws11 @{
    "E2" = 0.0522; "E3" = 0.0522; "E4" = 0.0522
    "O4" = 0; "P4" = 0; "Q4" = 0; "R4" = 0; "S4" = 0; "T4" = 0; "U4" = 0; "V4" = 0; "W4" = 0

    "E7" = 0.0214; "E8" = 0.0214; "E9" = 0.0214
    "O9" = 0.0024
    "P9" = 0.00305833333333333; "Q9" = 0.00305833333333333; "R9" = 0.009175
    "S9" = 0.00305833333333333; "T9" = 0.00305833333333333; "U9" = 0.00305833333333333
    "V9" = 0.009175; "W9" = 0.0367
}

# ---------------------------------------------------------------------------
# 3) Ciserano Italy
#    E2 and E3 (ytd) become blank; row 4 (Commit/Forecast) is removed entirely.
# ---------------------------------------------------------------------------
$ws12 = $wb.Worksheets.Item("Ciserano Italy")

$ws12.Range("E2").ClearContents()
$ws12.Range("E3").ClearContents()
$ws12.Rows.Item(4).Delete()

# ---------------------------------------------------------------------------
# 4) Changzhou Epc China
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Changzhou Epc China")

Set-Cells $ws3 @{
    "E2" = 0.0685; "E3" = 0.0685; "E4" = 0.0685
    "H4" = 0.0169; "J4" = 0.0169; "K4" = 0.0169; "M4" = 0.0351; "N4" = 0.0515
    "O4" = 0; "P4" = 0; "Q4" = 0; "R4" = 0; "S4" = 0; "T4" = 0; "U4" = 0; "V4" = 0; "W4" = 0

    "E8" = 0.0561; "E9" = 0.0561; "E10" = 0.0561
    "H10" = 0.0035; "I10" = 0.0104; "K10" = 0.0105; "L10" = 0.0178; "M10" = 0.0036; "N10" = 0.0318
    "O10" = 0.0036
    "P10" = 0.00801666666666667; "Q10" = 0.00801666666666667; "R10" = 0.02405
    "S10" = 0.00801666666666667; "T10" = 0.00801666666666667; "U10" = 0.00801666666666667
    "V10" = 0.02405; "W10" = 0.0962
}

# O7 (Internal Fill Rate / Commit-Forecast) goes from an explicit 0 to blank.
$ws3.Range("O7").ClearContents()

# ---------------------------------------------------------------------------
# 5) Jiaxing China
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Jiaxing China")

Set-Cells $ws4 @{
    "E2" = 0.048; "E3" = 0.048; "E4" = 0.048
    "O4" = 0; "P4" = 0; "Q4" = 0; "R4" = 0; "S4" = 0; "T4" = 0; "U4" = 0; "V4" = 0; "W4" = 0

    "E7" = 0.0317; "E8" = 0.0317; "E9" = 0.0317
    "O9" = 0.0071
    "P9" = 0.004525; "Q9" = 0.004525; "R9" = 0.013575
    "S9" = 0.004525; "T9" = 0.004525; "U9" = 0.004525
    "V9" = 0.013575; "W9" = 0.0543
}

# ---------------------------------------------------------------------------
# 6) Suzhou China
#    A new "Internal Fill Rate / Commit-Forecast" row is inserted as row 5;
#    the former Manufacturing Voluntary Turnover rows (5,6,7) shift down to
#    become rows 6,7,8, and several of their values are updated.
# ---------------------------------------------------------------------------
$ws9 = $wb.Worksheets.Item("Suzhou China")

$ws9.Rows.Item(5).Insert()

# New row 5: Internal Fill Rate / Commit/Forecast (all new content)
$ws9.Range("A5").Value = "PES"
$ws9.Range("B5").Value = "PES APAC"
$ws9.Range("C5").Value = "Suzhou China"
$ws9.Range("D5").Value = "Internal Fill Rate"
$ws9.Range("E5").Value = 0
$ws9.Range("F5").Value = "Commit/Forecast"
# G5:N5 remain blank (already blank after the row insert)
Set-Cells $ws9 @{
    "O5" = 0; "P5" = 0; "Q5" = 0; "R5" = 0; "S5" = 0; "T5" = 0; "U5" = 0; "V5" = 0; "W5" = 0
}

# Row 6 (previously row 5): Manufacturing Voluntary Turnover / PY Actual
Set-Cells $ws9 @{
    "E6" = 0.0604
    "K6" = 0.0036; "M6" = 0.0035; "N6" = 0.0071
    "O6" = 0.0141; "P6" = 0.007; "R6" = 0.021
    "S6" = 0.0035; "T6" = 0.0035; "V6" = 0.0069; "W6" = 0.0354
}

# Row 7 (previously row 6): Manufacturing Voluntary Turnover / AOP
Set-Cells $ws9 @{
    "E7" = 0.0604
    "K7" = 0.00324; "M7" = 0.00315; "N7" = 0.00639
    "O7" = 0.01269; "P7" = 0.0063; "R7" = 0.0189
    "S7" = 0.00315; "T7" = 0.00315; "V7" = 0.00621; "W7" = 0.03186
}

# Row 8 (previously row 7): Manufacturing Voluntary Turnover / Commit-Forecast
Set-Cells $ws9 @{
    "E8" = 0.0604
    "O8" = 0.0065
    "P8" = 0.008625; "Q8" = 0.008625; "R8" = 0.025875
    "S8" = 0.008625; "T8" = 0.008625; "U8" = 0.008625
    "V8" = 0.025875; "W8" = 0.1035
}
